$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.714.28'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '3.451.08'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.71'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.40'
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.450.10'
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("E9").Value = '  -5.26%  '
$ws.Range("E10").Value = '  +0.00%  '
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("E12").Value = '  -1.27%  '
$ws.Range("D13").Value = '4.041.75'
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("E15").Value = '  -2.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000174'
$ws.Range("E16").Value = '  -9.58%  '
$ws.Range("D17").Value = '64.740.35'
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").Value = '3.451.20'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.19'
$ws.Range("E19").Value = '  -3.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.80'
$ws.Range("E20").Value = '  -3.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.77'
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("E22").Value = '  -1.43%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.49'
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.532'
$ws.Range("E25").Value = '  -3.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000119'
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.994'
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.44'
$ws.Range("E30").Value = '  -3.40%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.07'
$ws.Range("E31").Value = '  -0.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.01'
$ws.Range("E32").Value = '  -1.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.25'
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.04'
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.59'
$ws.Range("E35").Value = '  -2.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.28'
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  -2.20%  '
$ws.Range("D38").Value = '2.880.33'
$ws.Range("E38").Value = '  -2.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0748'
$ws.Range("E39").Value = '  -3.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.40'
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.801'
$ws.Range("E41").Value = '  +3.50%  '
$ws.Range("E42").Value = '  -0.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.01'
$ws.Range("E43").Value = '  +0.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.51'
$ws.Range("E44").Value = '  -1.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.08'
$ws.Range("E45").Value = '  +1.99%  '
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.44'
$ws.Range("E47").Value = '  +13.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '322.15'
$ws.Range("E48").Value = '  +3.48%  '
$ws.Range("E49").Value = '  -1.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.46'
$ws.Range("E50").Value = '  -2.30%  '
$ws.Range("E51").Value = '  -1.89%  '
